$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 472
$ws.Range("F6").Value = 832
$ws.Range("F8").Value = 1224
$ws.Range("F12").Value = 704
$ws.Range("F14").Value = 519
$ws.Range("F18").Value = 2957
$ws.Range("F19").Value = 2630
$ws.Range("F24").Value = 234
$ws.Range("F26").Value = 5321
$ws.Range("F31").Value = 326

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1136
$ws.Range("F10").Value = 34
$ws.Range("F14").Value = 613
$ws.Range("F19").Value = 44
$ws.Range("F26").Value = 3957
$ws.Range("F30").Value = 200
$ws.Range("F31").Value = 53
$ws.Range("F34").Value = 34

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2475
$ws.Range("F6").Value = 1054
$ws.Range("F9").Value = 1337

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2475
$ws.Range("F6").Value = 1054
$ws.Range("F7").Value = 1337
$ws.Range("F11").Value = 472
$ws.Range("F12").Value = 832
$ws.Range("F14").Value = 1224
$ws.Range("F17").Value = 704
$ws.Range("F18").Value = 1136
$ws.Range("F19").Value = 1136
$ws.Range("F21").Value = 519
$ws.Range("F23").Value = 2957
$ws.Range("F24").Value = 2630
$ws.Range("F27").Value = 34
$ws.Range("F28").Value = 234
$ws.Range("F29").Value = 5321
$ws.Range("F32").Value = 613
$ws.Range("F33").Value = 613
$ws.Range("F37").Value = 326
$ws.Range("F40").Value = 44
$ws.Range("F47").Value = 200
$ws.Range("F48").Value = 53
